$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values between row 2 and row 3 for columns D, M, N, O, P, R, S

# Column D (Fecha)
$ws.Range("D2").Value = 44235
$ws.Range("D3").Value = 44417

# Column M (Volumen)
$ws.Range("M2").Value = 70
$ws.Range("M3").Value = 60

# Column N (Precio mínimo)
$ws.Range("N2").Value = 42000
$ws.Range("N3").Value = 26000

# Column O (Precio máximo)
$ws.Range("O2").Value = 42000
$ws.Range("O3").Value = 26000

# Column P (Precio promedio ponderado)
$ws.Range("P2").Value = 42000
$ws.Range("P3").Value = 26000

# Column R (Origen)
$ws.Range("R2").Value = "Región de Arica y Parinacota"
$ws.Range("R3").Value = "Perú"

# Column S (Precio $/Kg)
$ws.Range("S2").Value = 2333
$ws.Range("S3").Value = 1444
